{"js": "// The document contains three \"<id>...</id>\" markers, each built out of\n// three separate runs: \"<id>\" (Courier New / brown), the bare id text\n// (default font/black) and \"</id>\" (Courier New / brown). The edit\n// collapses each trio into a single Courier-New run and renumbers the\n// id from \"p015r_aN\" to \"p015r_N\".\nconst renames = [\n  [\"p015r_a1\", \"p015r_1\"],\n  [\"p015r_a2\", \"p015r_2\"],\n  [\"p015r_a3\", \"p015r_3\"],\n];\n\nfor (const [oldId, newId] of renames) {\n  const results = context.document.body.search(\"<id>\" + oldId + \"</id>\", {\n    matchCase: true,\n  });\n  results.load(\"text\");\n  await context.sync();\n\n  for (const r of results.items) {\n    // Replacing the whole \"<id>...</id>\" range with plain text merges the\n    // three runs into one, adopting the (Courier New) formatting of the\n    // range's leading run - matching the target markup.\n    r.insertText(\"<id>\" + newId + \"</id>\", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document has three \"<id>...</id>\" markers. Each one is currently\n# split across three runs - \"<id>\" (Courier New / brown), the bare id\n# text (default font/black) and \"</id>\" (Courier New / brown). Renumber\n# each id from \"p015r_aN\" to \"p015r_N\"; Find/Replace across the whole\n# tag merges the three runs into a single Courier-New run, matching the\n# target markup.\n$d = $word.ActiveDocument\n\n$renames = @(\n  @{Old = \"p015r_a1\"; New = \"p015r_1\"},\n  @{Old = \"p015r_a2\"; New = \"p015r_2\"},\n  @{Old = \"p015r_a3\"; New = \"p015r_3\"}\n)\n\nforeach ($pair in $renames) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = \"<id>\" + $pair.Old + \"</id>\"\n  $find.Replacement.Text = \"<id>\" + $pair.New + \"</id>\"\n\n  # wdFindContinue = 1, wdReplaceAll = 2\n  $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
